# Generate Report for Handoff
# The 39613f77-... file (and the 089a9d3e-... file's per-language Status)
# has moved from "Handed back: in sync with en-US" to "Ready for handoff".
# The handback for 39613f77 came back on a stale source version, so the
# handback timestamp advances and an Error Detail message is recorded on
# that row (both language sheets). The Overview roll-up reflects the
# 39613f77 row's new status/time, and the "Error Detail" column is widened
# so the long message is readable.

$wb = $excel.ActiveWorkbook

$readyForHandoff = "Ready for handoff"

$newHandbackZh = "2016-08-15 18:43:37"
$newHoDate     = "2016-08-15 18:43:42"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/52540e335c8d0077980edb9e8993495851a78b9a/e2e/39613f77-56ff-4866-a0b3-591de88e5561.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b09d9016be321243921aaf8ea879fac54bc5295b/e2e/39613f77-56ff-4866-a0b3-591de88e5561.md."

# ---- Overview sheet -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
# Row 3 = the 39613f77-... file
$overview.Range("E3").Value = $readyForHandoff   # zh-cn status
$overview.Range("F3").Value = $readyForHandoff   # de-de status
$overview.Range("G3").Value = $newHoDate         # Latest HO Xliff Generate Date

# ---- zh-cn sheet ------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $readyForHandoff        # 089a9d3e row Status
$zhcn.Range("C3").Value = $readyForHandoff        # 39613f77 row Status
$zhcn.Range("H3").Value = $newHandbackZh          # 39613f77 row Latest Handback DateTime
$zhcn.Range("P3").Value = $errorDetail            # 39613f77 row Error Detail
$zhcn.Range("P1").ColumnWidth = 39.16666667       # widen Error Detail column to 40

# ---- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $readyForHandoff        # 089a9d3e row Status
$dede.Range("C3").Value = $readyForHandoff        # 39613f77 row Status
$dede.Range("H3").Value = $newHoDate              # 39613f77 row Latest Handback DateTime
$dede.Range("P3").Value = $errorDetail            # 39613f77 row Error Detail
$dede.Range("P1").ColumnWidth = 39.16666667       # widen Error Detail column to 40
